$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data update: new "Patente/Motor/Chasis" values for the Regresión Preprod
# 06-05-2021 run (RGM011 -> RGM015).
$ws.Range("W2").Value = "RGM015"
$ws.Range("X2").Value = "1234567RGM015"
$ws.Range("Y2").Value = "1234567RGM015"

# Update the view: the sheet was left scrolled so column P is the left-most
# visible column, with Y3 as the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 16
$win.ScrollRow = 1
[void]$ws.Range("Y3").Select()
